$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1040.5454
$ws.Range("I38").Value = 546.9474
$ws.Range("K38").Value = 1640.8422
$ws.Range("M38").Value = -1268.8422
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 846.5294
$ws.Range("I2").Value = 492.92856
$ws.Range("K2").Value = 492.92856
$ws.Range("M2").Value = -379.92856
$ws.Range("H32").Value = 4412.273
$ws.Range("I32").Value = 4384.2856
$ws.Range("K32").Value = 4384.2856
$ws.Range("M32").Value = -4097.2856
$ws.Range("H44").Value = 52750
$ws.Range("J44").Value = 52750
$ws.Range("L44").Value = 52750
$ws.Range("N44").Value = -53726
$ws.Range("H61").Value = 4858.1875
$ws.Range("I61").Value = 1645.2858
$ws.Range("K61").Value = 1645.2858
$ws.Range("M61").Value = -1433.2858
$ws.Range("H74").Value = 1229.6
$ws.Range("I74").Value = 1229.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1229.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -355.5999999999999
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1229.6
$ws.Range("I77").Value = 1229.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6148
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1780
$ws.Range("N77").ClearContents()
$ws.Range("H116").Value = 846.5294
$ws.Range("I116").Value = 492.92856
$ws.Range("K116").Value = 492.92856
$ws.Range("M116").Value = 1801.07144
$ws.Range("H132").Value = 2221.2856
$ws.Range("I132").Value = 2091.5
$ws.Range("K132").Value = 6274.5
$ws.Range("M132").Value = -3744.5
$ws.Range("H136").Value = 4858.1875
$ws.Range("I136").Value = 1645.2858
$ws.Range("K136").Value = 4935.857400000001
$ws.Range("M136").Value = -2385.857400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 846.5294
$ws.Range("I3").Value = 492.92856
$ws.Range("K3").Value = 492.92856
$ws.Range("M3").Value = -378.92856
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 100
$ws.Range("M29").Value = 189
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1272.9231
$ws.Range("I58").Value = 1316.5
$ws.Range("J58").Value = 750
$ws.Range("K58").Value = 1316.5
$ws.Range("L58").Value = 750
$ws.Range("M58").Value = -1113.5
$ws.Range("N58").Value = -1156
$ws.Range("H100").Value = 100000
$ws.Range("J100").Value = 100000
$ws.Range("L100").Value = 100000
$ws.Range("N100").Value = -102164
$ws.Range("H122").Value = 3706.9092
$ws.Range("I122").Value = 3642.111
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 10926.333
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = -8476.332999999999
$ws.Range("N122").Value = -16895.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2848.5
$ws.Range("I132").Value = 2979.2
$ws.Range("J132").Value = 2195
$ws.Range("K132").Value = 8937.599999999999
$ws.Range("L132").Value = 6585
$ws.Range("M132").Value = -6407.599999999999
$ws.Range("N132").Value = -11645
$ws.Range("H136").Value = 1272.9231
$ws.Range("I136").Value = 1316.5
$ws.Range("J136").Value = 750
$ws.Range("K136").Value = 3949.5
$ws.Range("L136").Value = 2250
$ws.Range("M136").Value = -1399.5
$ws.Range("N136").Value = -7350
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 348.875
$ws.Range("I11").Value = 279
$ws.Range("K11").Value = 837
$ws.Range("M11").Value = -697
$ws.Range("H26").Value = 290.625
$ws.Range("I26").Value = 97.5
$ws.Range("K26").Value = 292.5
$ws.Range("M26").Value = -4.5
$ws.Range("H56").Value = 13697
$ws.Range("I56").Value = 13697
$ws.Range("K56").Value = 13697
$ws.Range("M56").Value = -13167
$ws.Range("H63").Value = 16364.75
$ws.Range("I63").Value = 16820
$ws.Range("J63").Value = 14999
$ws.Range("K63").Value = 50460
$ws.Range("L63").Value = 44997
$ws.Range("M63").Value = -49711
$ws.Range("N63").Value = -46495
$ws.Range("H66").Value = 16364.75
$ws.Range("I66").Value = 16820
$ws.Range("J66").Value = 14999
$ws.Range("K66").Value = 151380
$ws.Range("L66").Value = 134991
$ws.Range("M66").Value = -147636
$ws.Range("N66").Value = -142479
$ws.Range("I132").Value = 1049.5
$ws.Range("J132").Value = 15249.75
$ws.Range("K132").Value = 9445.5
$ws.Range("L132").Value = 137247.75
$ws.Range("M132").Value = -6915.5
$ws.Range("N132").Value = -142307.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2446.923
$ws.Range("I80").Value = 2099.5
$ws.Range("J80").Value = 2510.0908
$ws.Range("K80").Value = 2099.5
$ws.Range("L80").Value = 2510.0908
$ws.Range("M80").Value = -1101.5
$ws.Range("N80").Value = -4506.0908
$ws.Range("H83").Value = 2446.923
$ws.Range("I83").Value = 2099.5
$ws.Range("J83").Value = 2510.0908
$ws.Range("K83").Value = 10497.5
$ws.Range("L83").Value = 12550.454
$ws.Range("M83").Value = -5505.5
$ws.Range("N83").Value = -22534.454
$ws.Range("H126").Value = 3999.5
$ws.Range("I126").Value = 3666
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 10998
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -8528
$ws.Range("N126").Value = -19940
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H16").Value = 1829.5
$ws.Range("I16").Value = 1836.6
$ws.Range("J16").Value = 1794
$ws.Range("K16").Value = 1836.6
$ws.Range("L16").Value = 1794
$ws.Range("M16").Value = -1666.6
$ws.Range("N16").Value = -2134
$ws.Range("H22").Value = 2753.35
$ws.Range("I22").Value = 2155.5715
$ws.Range("K22").Value = 2155.5715
$ws.Range("M22").Value = -1860.5715
$ws.Range("H27").Value = 2753.35
$ws.Range("I27").Value = 2155.5715
$ws.Range("K27").Value = 2155.5715
$ws.Range("M27").Value = -2048.5715
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H46").Value = 3877.4
$ws.Range("I46").Value = 3795.6667
$ws.Range("K46").Value = 3795.6667
$ws.Range("M46").Value = -3607.6667
$ws.Range("H55").Value = 999
$ws.Range("I55").Value = 733.44446
$ws.Range("K55").Value = 733.44446
$ws.Range("M55").Value = -560.44446
$ws.Range("H68").Value = 2799.25
$ws.Range("J68").Value = 2799.25
$ws.Range("L68").Value = 2799.25
$ws.Range("N68").Value = -4297.25
$ws.Range("H71").Value = 2799.25
$ws.Range("J71").Value = 2799.25
$ws.Range("L71").Value = 13996.25
$ws.Range("N71").Value = -21484.25
$ws.Range("H93").Value = 1642.5
$ws.Range("I93").Value = 1642.5
$ws.Range("K93").Value = 1642.5
$ws.Range("M93").Value = -394.5
$ws.Range("H122").Value = 7726.636
$ws.Range("I122").Value = 7999
$ws.Range("J122").Value = 7666.1113
$ws.Range("K122").Value = 23997
$ws.Range("L122").Value = 22998.3339
$ws.Range("M122").Value = -21547
$ws.Range("N122").Value = -27898.3339
$ws.Range("H132").Value = 2949.1428
$ws.Range("J132").Value = 3935
$ws.Range("L132").Value = 11805
$ws.Range("N132").Value = -16865
$ws.Range("H136").Value = 4499.75
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1000001.5
$ws.Range("I8").Value = 1000001.5
$ws.Range("K8").Value = 1000001.5
$ws.Range("M8").Value = -999861.5
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 4
$ws.Range("K17").Value = 4
$ws.Range("M17").Value = 168
$ws.Range("H62").Value = 17779.176
$ws.Range("I62").Value = 28860
$ws.Range("J62").Value = 13162.167
$ws.Range("K62").Value = 28860
$ws.Range("L62").Value = 13162.167
$ws.Range("M62").Value = -28236
$ws.Range("N62").Value = -14410.167
$ws.Range("H65").Value = 17779.176
$ws.Range("I65").Value = 28860
$ws.Range("J65").Value = 13162.167
$ws.Range("K65").Value = 144300
$ws.Range("L65").Value = 65810.83499999999
$ws.Range("M65").Value = -141180
$ws.Range("N65").Value = -72050.83499999999
$ws.Range("H100").Value = 1315.7693
$ws.Range("I100").Value = 1300.4166
$ws.Range("K100").Value = 2600.8332
$ws.Range("M100").Value = -2059.8332
$ws.Range("H107").Value = 1394.5714
$ws.Range("I107").Value = 1394.5714
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4183.7142
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2263.7142
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 1017
$ws.Range("J113").Value = 562.3333
$ws.Range("L113").Value = 1686.9999
$ws.Range("N113").Value = -6026.9999
$ws.Range("H126").Value = 3523
$ws.Range("J126").Value = 5211
$ws.Range("L126").Value = 15633
$ws.Range("N126").Value = -20573
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 4272.364
$ws.Range("I132").Value = 4299.143
$ws.Range("K132").Value = 12897.429
$ws.Range("M132").Value = -10367.429
